# Apply the changes described by the diff:
#  1. On sheet "s1": fill column Z (rows 2-71) with the "answer" values
#     (1/2/3 marking the correct answer choice for each question row).
#  2. On sheet "s1": move the selection to P73 (scrolled so G43 area is
#     visible).
#  3. On sheet "Sheet1": select C3 so the saved view shows that cell
#     selected.

$wb = $excel.ActiveWorkbook

# --- Sheet "s1": Z column values (rows 2-71) ----------------------------
$ws1 = $wb.Worksheets.Item("s1")

$zValues = @(2,1,2,1,2,3,3,1,2,2,2,3,3,2,2,1,2,2,2,1,2,2,1,2,3,2,2,3,2,1,3,1,3,2,2,2,3,2,1,2,2,2,3,1,3,3,2,1,1,3,1,3,2,3,2,1,3,2,1,3,3,2,2,3,1,1,2,3,3,2)

for ($i = 0; $i -lt $zValues.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 26).Value = $zValues[$i]
}

# --- Sheet "s1": scroll/selection update --------------------------------
$ws1.Activate()
# Bring the G43 area into view, then land the active selection on P73
# (matches the recorded sheetView/selection in the saved workbook).
$ws1.Range("G43").Select()
$ws1.Range("P73").Select()

# --- Sheet "Sheet1": select C3 ------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Activate()
$ws2.Range("C3").Select()
